$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.909.07"
$ws.Range("E2").Value = "  +2.86%  "

$ws.Range("D3").Value = "2.421.73"
$ws.Range("E3").Value = "  +2.58%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.32%  "

$ws.Range("D5").Value = "'552.06"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").Value = "'137.88"
$ws.Range("E6").Value = "  +3.47%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.23%  "

$ws.Range("E8").Value = "  +3.38%  "

$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("E10").Value = "  +0.14%  "

$ws.Range("B11").Value = "TRON"
$ws.Range("C11").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D11").Value = "'0.148"
$ws.Range("E11").Value = "  -1.93%  "

$ws.Range("B12").Value = "Cardano"
$ws.Range("C12").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D12").Value = "'0.355"
$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("D13").Value = "'25.19"
$ws.Range("E13").Value = "  +4.11%  "

$ws.Range("D14").Value = "2.853.96"
$ws.Range("E14").Value = "  +2.96%  "

$ws.Range("D15").Value = "59.858.58"
$ws.Range("E15").Value = "  +3.27%  "

$ws.Range("E16").Value = "  +0.82%  "

$ws.Range("D17").Value = "2.394.25"
$ws.Range("E17").Value = "  +0.54%  "

$ws.Range("E18").Value = "  +2.35%  "

$ws.Range("E19").Value = "  +1.39%  "

$ws.Range("D20").Value = "'330.36"
$ws.Range("E20").Value = "  -0.42%  "

$ws.Range("D21").Value = "'6.69"
$ws.Range("E21").Value = "  -3.39%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.22%  "

$ws.Range("D23").Value = "'66.14"
$ws.Range("E23").Value = "  +3.77%  "

$ws.Range("D24").Value = "'0.172"
$ws.Range("E24").Value = "  +0.99%  "

$ws.Range("D25").Value = "'8.84"
$ws.Range("E25").Value = "  +6.49%  "

$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("E27").Value = "  +3.85%  "

$ws.Range("D28").Value = "0.0₃0777"
$ws.Range("E28").Value = "  +4.79%  "

$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("D30").Value = "'170.28"
$ws.Range("E30").Value = "  -0.23%  "

$ws.Range("E31").Value = "  -0.83%  "

$ws.Range("D32").Value = "'18.67"
$ws.Range("E32").Value = "  +1.34%  "

$ws.Range("E33").Value = "  +1.35%  "

$ws.Range("D35").Value = "'1.29"
$ws.Range("E35").Value = "  +4.23%  "

$ws.Range("E36").Value = "  +0.07%  "

$ws.Range("D37").Value = "'4.21"
$ws.Range("E37").Value = "  +0.92%  "

$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").Value = "'39.60"
$ws.Range("E39").Value = "  -1.83%  "

$ws.Range("D40").Value = "'0.410"
$ws.Range("E40").Value = "  -3.99%  "

$ws.Range("D41").Value = "'313.23"
$ws.Range("E41").Value = "  +8.63%  "

$ws.Range("D42").Value = "'3.68"
$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").Value = "'138.99"
$ws.Range("E43").Value = "  -1.36%  "

$ws.Range("D44").Value = "'0.0971"
$ws.Range("E44").Value = "  +1.25%  "

$ws.Range("E45").Value = "  +0.65%  "

$ws.Range("D46").Value = "'19.46"
$ws.Range("E46").Value = "  +4.40%  "

$ws.Range("D47").Value = "'0.579"
$ws.Range("E47").Value = "  +2.19%  "

$ws.Range("E48").Value = "  +0.94%  "

$ws.Range("E49").Value = "  -2.18%  "

$ws.Range("D50").Value = "'17.62"
$ws.Range("E50").Value = "  +1.11%  "

$ws.Range("E51").Value = "  +0.33%  "
